$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work plan")

# --- Row 42: update End date, Work results, Next steps, Notes ---
$ws.Cells.Item(42,4).Style = "Geras"
$ws.Cells.Item(42,4).Value = "16/03/2020"

$ws.Cells.Item(42,6).Value = "INSTALL.md file explains how to run this software"
$ws.Cells.Item(42,7).Value = "Start writing final report"

$ws.Cells.Item(42,8).WrapText = $true
$ws.Cells.Item(42,8).Value = "Didn’t know how much in detail I need to explain. Depends on the audience. Assumed that reader is semi skilled"

$ws.Rows.Item(42).RowHeight = 43.2

# --- Row 50 (new log entry) ---
$ws.Cells.Item(50,5).WrapText = $true
$ws.Cells.Item(50,5).Value = "Color support/ including color schema for blind people"

$ws.Cells.Item(50,2).NumberFormat = "@"
$ws.Cells.Item(50,2).Value = "20/03/2020"

$ws.Cells.Item(50,3).NumberFormat = "@"
$ws.Cells.Item(50,3).Value = "20/03/2020"

$ws.Cells.Item(50,4).Style = "Geras"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "20/03/2020"

$ws.Cells.Item(50,6).Value = "Code is working and pushed on github"
$ws.Cells.Item(50,7).Value = "Mention about it in the report"

# --- Row 45: update End date, Work results, Next steps, Notes ---
$ws.Cells.Item(45,6).WrapText = $true
$ws.Cells.Item(45,6).Value = "Completely rewrote the specifications, eliminated I, we etc. SMART target list is still missing"

$ws.Cells.Item(45,7).Value = "Complete the SMART list"

$ws.Cells.Item(45,8).WrapText = $true
$ws.Cells.Item(45,8).Value = "Don’t know how to write SMART list, need lecturer assistance"

$ws.Cells.Item(45,4).Style = "Geras"
$ws.Cells.Item(45,4).Value = "20/03/2020"

$ws.Rows.Item(45).RowHeight = 28.8

# --- Row 46: update End date, Work results ---
$ws.Cells.Item(46,6).WrapText = $true
$ws.Cells.Item(46,6).Value = "Have a paper note with the topics to write about. Also read the marking criteria"

$ws.Cells.Item(46,4).Style = "Geras"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "20/03/2020"

$ws.Rows.Item(46).RowHeight = 28.8

# --- View: select F50 ---
$ws.Range("F50").Select() | Out-Null

Write-Host "done"
